$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nueva fila de inventario (fila 11)
$ws.Range("A11").Value = "10ENFV"
$ws.Range("B11").Value = "Almohadilla + Chip Epson C9344"
$ws.Range("C11").Value = "WF 2810 2830 2835 2850 2851 2930 2950, XP 2100 2105 3100 3105 4100 4101 4105, L3550 L3560 L3590 L5550 L5560 L5590"
$ws.Range("D11").Value = 35000
$ws.Range("E11").Value = 200000
$ws.Range("F11").Value = 11
$ws.Range("G11").Value = 9
$ws.Range("H11").Formula = "=(E11-D11)*G11"
$ws.Range("I11").Formula = "=D11*F11"
$ws.Range("J11").Value = 385000
